$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.181.50"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "1.868.79"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "'336.94"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("D7").Value = "'0.4712"
$ws.Range("E7").Value = "  +1.58%  "
$ws.Range("D8").Value = "'0.3925"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").Value = "'46.85"
$ws.Range("E9").Value = "  +2.09%  "
$ws.Range("D10").Value = "'0.07975"
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "'1.004"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "'21.74"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.864.69"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.982"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.273"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "'0.06596"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "'17.76"
$ws.Range("E20").Value = "  +3.94%  "
$ws.Range("D21").Value = "'0.9990"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").Value = "28.177.27"
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").Value = "'5.440"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").Value = "'11.05"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").Value = "'2.294"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").Value = "2.079.98"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "'159.12"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").Value = "'19.81"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").Value = "'5.489"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").Value = "'119.72"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "'0.9739"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").Value = "'3.573"
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D35").Value = "'1.381"
$ws.Range("E35").Value = "  +3.48%  "
$ws.Range("D36").Value = "'5.342"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").Value = "'0.02272"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("D38").Value = "'0.06097"
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").Value = "'8.405"
$ws.Range("E39").Value = "  +1.32%  "
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "'0.5983"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("D42").Value = "'0.9992"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").Value = "'0.1884"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("D46").Value = "'0.5629"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").Value = "'12.14"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").Value = "'1.979"
$ws.Range("E48").Value = "  +3.90%  "
$ws.Range("D49").Value = "'0.06857"
$ws.Range("E49").Value = "  +2.48%  "
$ws.Range("D50").Value = "'111.56"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").Value = "'1.055"
$ws.Range("E51").Value = "  +0.56%  "
